$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "time_taken" header in F1, matching the style of the other header cells
$ws.Cells.Item(1, 6).Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the time_taken column (F2:F45) with the recorded timestamps as plain text
$ws.Cells.Item(2, 6).Value = "2021-10-05 13:39:46.481651"
$ws.Cells.Item(3, 6).Value = "2021-10-05 13:39:46.481662"
$ws.Cells.Item(4, 6).Value = "2021-10-05 13:39:46.481665"
$ws.Cells.Item(5, 6).Value = "2021-10-05 13:39:46.481668"
$ws.Cells.Item(6, 6).Value = "2021-10-05 13:39:46.481671"
$ws.Cells.Item(7, 6).Value = "2021-10-05 13:39:46.481673"
$ws.Cells.Item(8, 6).Value = "2021-10-05 13:39:46.481676"
$ws.Cells.Item(9, 6).Value = "2021-10-05 13:39:46.481678"
$ws.Cells.Item(10, 6).Value = "2021-10-05 13:39:46.481681"
$ws.Cells.Item(11, 6).Value = "2021-10-05 13:39:46.481683"
$ws.Cells.Item(12, 6).Value = "2021-10-05 13:39:46.481686"
$ws.Cells.Item(13, 6).Value = "2021-10-05 13:39:46.481688"
$ws.Cells.Item(14, 6).Value = "2021-10-05 13:39:46.481691"
$ws.Cells.Item(15, 6).Value = "2021-10-05 13:39:46.481693"
$ws.Cells.Item(16, 6).Value = "2021-10-05 13:39:46.481696"
$ws.Cells.Item(17, 6).Value = "2021-10-05 13:39:46.481698"
$ws.Cells.Item(18, 6).Value = "2021-10-05 13:39:46.481701"
$ws.Cells.Item(19, 6).Value = "2021-10-05 13:39:46.481704"
$ws.Cells.Item(20, 6).Value = "2021-10-05 13:39:46.481706"
$ws.Cells.Item(21, 6).Value = "2021-10-05 13:39:46.481709"
$ws.Cells.Item(22, 6).Value = "2021-10-05 13:39:46.481711"
$ws.Cells.Item(23, 6).Value = "2021-10-05 13:39:46.481714"
$ws.Cells.Item(24, 6).Value = "2021-10-05 13:39:46.481716"
$ws.Cells.Item(25, 6).Value = "2021-10-05 13:39:46.481719"
$ws.Cells.Item(26, 6).Value = "2021-10-05 13:39:46.481721"
$ws.Cells.Item(27, 6).Value = "2021-10-05 13:39:46.481724"
$ws.Cells.Item(28, 6).Value = "2021-10-05 13:39:46.481726"
$ws.Cells.Item(29, 6).Value = "2021-10-05 13:39:46.481729"
$ws.Cells.Item(30, 6).Value = "2021-10-05 13:39:46.481731"
$ws.Cells.Item(31, 6).Value = "2021-10-05 13:39:46.481734"
$ws.Cells.Item(32, 6).Value = "2021-10-05 13:39:46.481736"
$ws.Cells.Item(33, 6).Value = "2021-10-05 13:39:46.481739"
$ws.Cells.Item(34, 6).Value = "2021-10-05 13:39:46.481741"
$ws.Cells.Item(35, 6).Value = "2021-10-05 13:39:46.481744"
$ws.Cells.Item(36, 6).Value = "2021-10-05 13:39:46.481746"
$ws.Cells.Item(37, 6).Value = "2021-10-05 13:39:46.481749"
$ws.Cells.Item(38, 6).Value = "2021-10-05 13:39:46.481751"
$ws.Cells.Item(39, 6).Value = "2021-10-05 13:39:46.481753"
$ws.Cells.Item(40, 6).Value = "2021-10-05 13:39:46.481756"
$ws.Cells.Item(41, 6).Value = "2021-10-05 13:39:46.481758"
$ws.Cells.Item(42, 6).Value = "2021-10-05 13:39:46.481761"
$ws.Cells.Item(43, 6).Value = "2021-10-05 13:39:46.481764"
$ws.Cells.Item(44, 6).Value = "2021-10-05 13:39:46.481766"
$ws.Cells.Item(45, 6).Value = "2021-10-05 13:39:46.481769"
